$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(76, 8).Value = 3911669.8
$ws.Cells.Item(76, 9).Value = 5862606
$ws.Cells.Item(76, 11).Value = 5862606
$ws.Cells.Item(76, 13).Value = -5862291
$ws.Cells.Item(79, 8).Value = 3911669.8
$ws.Cells.Item(79, 9).Value = 5862606
$ws.Cells.Item(79, 11).Value = 5862606
$ws.Cells.Item(79, 13).Value = -5861514
$ws.Cells.Item(80, 8).Value = 911.8333
$ws.Cells.Item(80, 9).Value = 868.6
$ws.Cells.Item(80, 11).Value = 2605.8
$ws.Cells.Item(80, 13).Value = -1607.8
$ws.Cells.Item(83, 8).Value = 911.8333
$ws.Cells.Item(83, 9).Value = 868.6
$ws.Cells.Item(83, 11).Value = 7817.400000000001
$ws.Cells.Item(83, 13).Value = -2825.400000000001
$ws.Cells.Item(125, 8).Value = 1283
$ws.Cells.Item(125, 9).Value = 1650
$ws.Cells.Item(125, 10).Value = 1038.3334
$ws.Cells.Item(125, 11).Value = 14850
$ws.Cells.Item(125, 12).Value = 9345.000599999999
$ws.Cells.Item(125, 13).Value = -12390
$ws.Cells.Item(125, 14).Value = -14265.0006
$ws.Cells.Item(137, 8).Value = 1904.5
$ws.Cells.Item(137, 9).Value = 1635.6364
$ws.Cells.Item(137, 10).Value = 2233.111
$ws.Cells.Item(137, 11).Value = 4906.9092
$ws.Cells.Item(137, 12).Value = 6699.333
$ws.Cells.Item(137, 13).Value = -2356.9092
$ws.Cells.Item(137, 14).Value = -11799.333
$ws.Cells.Item(138, 8).Value = 2535.3518
$ws.Cells.Item(138, 9).Value = 3249.1052
$ws.Cells.Item(138, 10).Value = 2147.8857
$ws.Cells.Item(138, 11).Value = 9747.3156
$ws.Cells.Item(138, 12).Value = 6443.657099999999
$ws.Cells.Item(138, 13).Value = -4607.3156
$ws.Cells.Item(138, 14).Value = -16723.6571
$ws.Cells.Item(141, 8).Value = 3123.3076
$ws.Cells.Item(141, 9).Value = 2418.5
$ws.Cells.Item(141, 11).Value = 7255.5
$ws.Cells.Item(141, 13).Value = -2075.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1172.875
$ws.Cells.Item(2, 9).Value = 1226.6
$ws.Cells.Item(2, 11).Value = 1226.6
$ws.Cells.Item(2, 13).Value = -1113.6
$ws.Cells.Item(32, 8).Value = 2950.7778
$ws.Cells.Item(32, 9).Value = 2058.638
$ws.Cells.Item(32, 10).Value = 6646.7856
$ws.Cells.Item(32, 11).Value = 2058.638
$ws.Cells.Item(32, 12).Value = 6646.7856
$ws.Cells.Item(32, 13).Value = -1771.638
$ws.Cells.Item(32, 14).Value = -7220.7856
$ws.Cells.Item(45, 8).Value = 1346.7333
$ws.Cells.Item(45, 9).Value = 859.7143
$ws.Cells.Item(45, 11).Value = 859.7143
$ws.Cells.Item(45, 13).Value = -482.7143
$ws.Cells.Item(61, 8).Value = 2493.3914
$ws.Cells.Item(61, 9).Value = 1317.3125
$ws.Cells.Item(61, 11).Value = 1317.3125
$ws.Cells.Item(61, 13).Value = -1105.3125
$ws.Cells.Item(74, 8).Value = 1800.7693
$ws.Cells.Item(74, 9).Value = 1749.8334
$ws.Cells.Item(74, 10).Value = 1844.4286
$ws.Cells.Item(74, 11).Value = 1749.8334
$ws.Cells.Item(74, 12).Value = 1844.4286
$ws.Cells.Item(74, 13).Value = -875.8334
$ws.Cells.Item(74, 14).Value = -3592.4286
$ws.Cells.Item(77, 8).Value = 1800.7693
$ws.Cells.Item(77, 9).Value = 1749.8334
$ws.Cells.Item(77, 10).Value = 1844.4286
$ws.Cells.Item(77, 11).Value = 8749.166999999999
$ws.Cells.Item(77, 12).Value = 9222.143
$ws.Cells.Item(77, 13).Value = -4381.166999999999
$ws.Cells.Item(77, 14).Value = -17958.143
$ws.Cells.Item(116, 8).Value = 1172.875
$ws.Cells.Item(116, 9).Value = 1226.6
$ws.Cells.Item(116, 11).Value = 1226.6
$ws.Cells.Item(116, 13).Value = 1067.4
$ws.Cells.Item(122, 8).Value = 638.2
$ws.Cells.Item(122, 9).Value = 638.2
$ws.Cells.Item(122, 11).Value = 1914.6
$ws.Cells.Item(122, 13).Value = 535.3999999999999
$ws.Cells.Item(132, 8).Value = 1384.0656
$ws.Cells.Item(132, 9).Value = 1081.4286
$ws.Cells.Item(132, 10).Value = 2619.8333
$ws.Cells.Item(132, 11).Value = 3244.2858
$ws.Cells.Item(132, 12).Value = 7859.499899999999
$ws.Cells.Item(132, 13).Value = -714.2857999999997
$ws.Cells.Item(132, 14).Value = -12919.4999
$ws.Cells.Item(136, 8).Value = 2493.3914
$ws.Cells.Item(136, 9).Value = 1317.3125
$ws.Cells.Item(136, 11).Value = 3951.9375
$ws.Cells.Item(136, 13).Value = -1401.9375

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1172.875
$ws.Cells.Item(3, 9).Value = 1226.6
$ws.Cells.Item(3, 11).Value = 1226.6
$ws.Cells.Item(3, 13).Value = -1112.6
$ws.Cells.Item(99, 8).Value = 1549.3572
$ws.Cells.Item(99, 9).Value = 1243.4445
$ws.Cells.Item(99, 11).Value = 1243.4445
$ws.Cells.Item(99, 13).Value = 254.5554999999999
$ws.Cells.Item(134, 8).Value = 5948.643
$ws.Cells.Item(134, 9).Value = 6811.391
$ws.Cells.Item(134, 10).Value = 1980
$ws.Cells.Item(134, 11).Value = 20434.173
$ws.Cells.Item(134, 12).Value = 5940
$ws.Cells.Item(134, 13).Value = -17899.173
$ws.Cells.Item(134, 14).Value = -11010

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1803.9
$ws.Cells.Item(31, 9).Value = 1434.4286
$ws.Cells.Item(31, 10).Value = 2666
$ws.Cells.Item(31, 11).Value = 1434.4286
$ws.Cells.Item(31, 12).Value = 2666
$ws.Cells.Item(31, 13).Value = -1139.4286
$ws.Cells.Item(31, 14).Value = -3256
$ws.Cells.Item(34, 8).Value = 1803.9
$ws.Cells.Item(34, 9).Value = 1434.4286
$ws.Cells.Item(34, 10).Value = 2666
$ws.Cells.Item(34, 11).Value = 1434.4286
$ws.Cells.Item(34, 12).Value = 2666
$ws.Cells.Item(34, 13).Value = -1232.4286
$ws.Cells.Item(34, 14).Value = -3070

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(33, 8).Value = 93.8
$ws.Cells.Item(33, 9).Value = 55
$ws.Cells.Item(33, 10).Value = 103.5
$ws.Cells.Item(33, 11).Value = 330
$ws.Cells.Item(33, 12).Value = 621
$ws.Cells.Item(33, 13).Value = -47
$ws.Cells.Item(33, 14).Value = -1187
$ws.Cells.Item(44, 8).Value = 249.66667
$ws.Cells.Item(44, 9).Value = 150
$ws.Cells.Item(44, 11).Value = 450
$ws.Cells.Item(44, 13).Value = -52
$ws.Cells.Item(122, 8).Value = 1120.909
$ws.Cells.Item(122, 9).Value = 899.6667
$ws.Cells.Item(122, 10).Value = 1203.875
$ws.Cells.Item(122, 11).Value = 8097.0003
$ws.Cells.Item(122, 12).Value = 10834.875
$ws.Cells.Item(122, 13).Value = -5647.0003
$ws.Cells.Item(122, 14).Value = -15734.875
$ws.Cells.Item(131, 8).Value = 1674.4142
$ws.Cells.Item(131, 10).Value = 1743.1183
$ws.Cells.Item(131, 12).Value = 5229.3549
$ws.Cells.Item(131, 14).Value = -15309.3549
$ws.Cells.Item(140, 8).Value = 1792.4722
$ws.Cells.Item(140, 9).Value = 1010.05554
$ws.Cells.Item(140, 10).Value = 2574.889
$ws.Cells.Item(140, 11).Value = 3030.16662
$ws.Cells.Item(140, 12).Value = 7724.667
$ws.Cells.Item(140, 13).Value = 2149.83338
$ws.Cells.Item(140, 14).Value = -18084.667

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 4999.8
$ws.Cells.Item(70, 10).Value = 4499.75
$ws.Cells.Item(70, 12).Value = 4499.75
$ws.Cells.Item(70, 14).Value = -5039.75
$ws.Cells.Item(73, 8).Value = 4999.8
$ws.Cells.Item(73, 10).Value = 4499.75
$ws.Cells.Item(73, 12).Value = 4499.75
$ws.Cells.Item(73, 14).Value = -6371.75
$ws.Cells.Item(97, 8).Value = 1695.8334
$ws.Cells.Item(97, 9).Value = 1654.875
$ws.Cells.Item(97, 10).Value = 1777.75
$ws.Cells.Item(97, 11).Value = 1654.875
$ws.Cells.Item(97, 12).Value = 1777.75
$ws.Cells.Item(97, 13).Value = -1158.875
$ws.Cells.Item(97, 14).Value = -2769.75
$ws.Cells.Item(122, 8).Value = 2072.8
$ws.Cells.Item(122, 9).Value = 1527.7142
$ws.Cells.Item(122, 10).Value = 2549.75
$ws.Cells.Item(122, 11).Value = 4583.142599999999
$ws.Cells.Item(122, 12).Value = 7649.25
$ws.Cells.Item(122, 13).Value = -2133.142599999999
$ws.Cells.Item(122, 14).Value = -12549.25
$ws.Cells.Item(132, 8).Value = 2853.4243
$ws.Cells.Item(132, 9).Value = 2597.7727
$ws.Cells.Item(132, 11).Value = 7793.3181
$ws.Cells.Item(132, 13).Value = -5263.3181

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 3280.2144
$ws.Cells.Item(7, 9).Value = 2619.25
$ws.Cells.Item(7, 11).Value = 2619.25
$ws.Cells.Item(7, 13).Value = -2507.25
$ws.Cells.Item(16, 8).Value = 5697
$ws.Cells.Item(16, 9).Value = 15999.5
$ws.Cells.Item(16, 11).Value = 15999.5
$ws.Cells.Item(16, 13).Value = -15829.5
$ws.Cells.Item(68, 8).Value = 2302
$ws.Cells.Item(68, 10).Value = 5000
$ws.Cells.Item(68, 12).Value = 5000
$ws.Cells.Item(68, 14).Value = -6498
$ws.Cells.Item(71, 8).Value = 2302
$ws.Cells.Item(71, 10).Value = 5000
$ws.Cells.Item(71, 12).Value = 25000
$ws.Cells.Item(71, 14).Value = -32488
$ws.Cells.Item(122, 8).Value = 5503.923
$ws.Cells.Item(122, 9).Value = 4443.4287
$ws.Cells.Item(122, 11).Value = 13330.2861
$ws.Cells.Item(122, 13).Value = -10880.2861
$ws.Cells.Item(126, 8).Value = 3280.2144
$ws.Cells.Item(126, 9).Value = 2619.25
$ws.Cells.Item(126, 11).Value = 7857.75
$ws.Cells.Item(126, 13).Value = -5387.75
$ws.Cells.Item(132, 8).Value = 2637.3794
$ws.Cells.Item(132, 9).Value = 2372.7
$ws.Cells.Item(132, 11).Value = 7118.099999999999
$ws.Cells.Item(132, 13).Value = -4588.099999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 5750.25
$ws.Cells.Item(126, 9).Value = 3200.4
$ws.Cells.Item(126, 10).Value = 10000
$ws.Cells.Item(126, 11).Value = 9601.200000000001
$ws.Cells.Item(126, 12).Value = 30000
$ws.Cells.Item(126, 13).Value = -7131.200000000001
$ws.Cells.Item(126, 14).Value = -34940
$ws.Cells.Item(136, 8).Value = 2558.5217
$ws.Cells.Item(136, 9).Value = 2096.9412
$ws.Cells.Item(136, 10).Value = 3866.3333
$ws.Cells.Item(136, 11).Value = 6290.823600000001
$ws.Cells.Item(136, 12).Value = 11598.9999
$ws.Cells.Item(136, 13).Value = -3740.823600000001
$ws.Cells.Item(136, 14).Value = -16698.9999
